# "show learnin category icon by image path"
#
# The "Имя файла иконки" (icon file name) column is repurposed to hold a
# URL ("Путь к картинке иконки" / "path to icon picture") that is turned
# into a clickable hyperlink, so the icon can be shown directly from its
# image path.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$iconUrl = "https://upload.wikimedia.org/wikipedia/commons/thumb/d/d6/Square_Flag_of_the_United_Kingdom.svg/1024px-Square_Flag_of_the_United_Kingdom.svg.png"

# Rename the header from "Имя файла иконки" to "Путь к картинке иконки".
$ws.Range("D1").Value = "Путь к картинке иконки"

# Turn the sample row's icon cell into a hyperlink pointing at the image,
# then show the URL itself as the cell text (this also applies Excel's
# built-in "Hyperlink" cell style automatically, same as the source diff).
$ws.Hyperlinks.Add($ws.Range("D2"), $iconUrl) | Out-Null
$ws.Range("D2").Value = $iconUrl

# Widen column D so the long URL/path is more readable.
$ws.Columns("D").ColumnWidth = 24

# Match the saved selection/active cell.
$ws.Range("D2").Select() | Out-Null
